$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FEB-2021")
$src = $wb.Worksheets.Item("OCT-2020")

# --- Row 14-15: "Holiday" block, mirrors the existing C7:G8 holiday block ---
$ws.Range("C14:G15").Merge()
$ws.Range("C7:G8").Copy()
$ws.Range("C14:G15").PasteSpecial(-4122)
$ws.Range("C14").Value = "Holiday"

# --- Row 16: QMVAR 2.0 / Adding new change request / WIP (same pattern as row 13) ---
$ws.Range("C13:F13").Copy()
$ws.Range("C16:F16").PasteSpecial(-4122)
$ws.Range("E9").Copy()
$ws.Range("G16").PasteSpecial(-4122)
$ws.Range("C16").Value = "QMVAR 2.0"
$ws.Range("D16").Value = "Adding new change request"
$ws.Range("F16").Value = "WIP"
$ws.Rows.Item(16).RowHeight = 28.8

# --- Row 17: nMVAR/QMVAR 2.0 combo (blank value set last, after rows 18/19) ---
$src.Range("C2").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("D9:D9").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("E9").Copy()
$ws.Range("E17:G17").PasteSpecial(-4122)
$ws.Range("E17:G17").WrapText = $true
$ws.Range("C17").Value = "1.nMVAR               2.QMVAR 2.0"
$ws.Range("D17").Value = "1. nMVAR issue fixing  2.Qmvar Issue fixing"
$ws.Rows.Item(17).RowHeight = 28.8

# --- Row 18: nMVAR / nMVAR issue fixing + WIP ---
$src.Range("C2").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("D9").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E17:G17").Copy()
$ws.Range("E18:G18").PasteSpecial(-4122)
$ws.Range("C18").Value = "1.nMVAR             "
$ws.Range("D18").Value = "1. nMVAR issue fixing  "
$ws.Range("F18").Value = "1. WIP  "

# --- Row 19: same as row 18 ---
$src.Range("C2").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$ws.Range("D9").Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("E17:G17").Copy()
$ws.Range("E19:G19").PasteSpecial(-4122)
$ws.Range("C19").Value = "1.nMVAR             "
$ws.Range("D19").Value = "1. nMVAR issue fixing  "
$ws.Range("F19").Value = "1. WIP  "

# --- back to Row 17: fill in F17 last so the new shared string is appended last ---
$ws.Range("F17").Value = "1.WIP               2.WIP"

# --- Row 20: same content/style as row 16 ---
$ws.Range("C16:G16").Copy()
$ws.Range("C20:G20").PasteSpecial(-4122)
$ws.Range("C20").Value = "QMVAR 2.0"
$ws.Range("D20").Value = "Adding new change request"
$ws.Range("F20").Value = "WIP"
$ws.Rows.Item(20).RowHeight = 28.8

# --- Selection moves to L16 ---
$ws.Activate()
$ws.Range("L16").Select()
